$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F6").Value = 268
$ws.Range("F7").Value = 13124
$ws.Range("F11").Value = 3870
$ws.Range("F12").Value = 6648
$ws.Range("F15").Value = 3498
$ws.Range("F16").Value = 39
$ws.Range("F20").Value = 75
$ws.Range("F22").Value = 3650
$ws.Range("F23").Value = 102
$ws.Range("F25").Value = 3522
$ws.Range("F26").Value = 3523
$ws.Range("F27").Value = 420
$ws.Range("F28").Value = 1912
$ws.Range("F30").Value = 236
$ws.Range("F31").Value = 6824
$ws.Range("F34").Value = 1596
$ws.Range("F35").Value = 2022
$ws.Range("F37").Value = 110
$ws.Range("F38").Value = 1072
$ws.Range("F40").Value = 221
$ws.Range("F43").Value = 1153
$ws.Range("F44").Value = 1148
$ws.Range("F47").Value = 1225
$ws.Range("F48").Value = 1816

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F11").Value = 12

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 461
$ws.Range("F3").Value = 628
$ws.Range("F4").Value = 28

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F6").Value = 461
$ws.Range("F7").Value = 628
$ws.Range("F8").Value = 28
$ws.Range("F9").Value = 268
$ws.Range("F10").Value = 13124
$ws.Range("F15").Value = 3871
$ws.Range("F16").Value = 6648
$ws.Range("F18").Value = 3498
$ws.Range("F19").Value = 39
$ws.Range("F22").Value = 75
$ws.Range("F28").Value = 3523
$ws.Range("F29").Value = 420
$ws.Range("F31").Value = 236
$ws.Range("F32").Value = 6824
$ws.Range("F36").Value = 1596
$ws.Range("F37").Value = 2022
$ws.Range("F39").Value = 110
$ws.Range("F40").Value = 1073
$ws.Range("F41").Value = 221
$ws.Range("F43").Value = 1153
$ws.Range("F44").Value = 1148
$ws.Range("F47").Value = 1816
